$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("R2M")

# Update the "WSA" row: version bumped from V4 to V5
$ws.Range("B3").Value = "V5"

# Update the "Additional" row: date bumped from 20200428 to 20200429
$ws.Range("B4").Value = 20200429

# Move the active selection to B3 (was B5)
$ws.Activate()
$ws.Range("B3").Select()
